# Applies the WSL_note.docx edit:
#  1. Strips the paragraph-mark rPr (<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>)
#     from the first paragraph ("Windows 進入 Linux 系統").
#  2. Appends, after the "wsl -d Ubuntu-22.04" paragraph:
#       - a blank paragraph
#       - a paragraph "Ctrl + O 儲存"
#       - a paragraph "Ctrl + X 退出" (with paragraph-mark rPr rFonts hint=eastAsia)

$d = $word.ActiveDocument
$pkgOpen = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1. Rewrite the first paragraph, dropping its <w:pPr> block -------------
$p1 = $d.Paragraphs(1)
$p1Body = '<w:body><w:p>' +
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Windows</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>進入</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Linux</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>系統</w:t></w:r>' +
  '</w:p></w:body>'
$p1.Range.InsertXML($pkgOpen + $p1Body + $pkgClose)

# --- 2. Append the three new paragraphs after the last paragraph ------------
$newBody = '<w:body>' +
  '<w:p/>' +
  '<w:p>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Ctrl + O </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>儲存</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Ctrl + X </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>退出</w:t></w:r>' +
  '</w:p>' +
  '</w:body>'

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($pkgOpen + $newBody + $pkgClose)
